$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-12 Wednesday" "2025-02-13 Thursday"

Replace-Text "964×4=3856" "915×4=3660"
Replace-Text "629×5=3145" "786×8=6288"
Replace-Text "208×4=832" "452×8=3616"
Replace-Text "297×2=594" "705×4=2820"
Replace-Text "159×7=1113" "181×4=724"

Replace-Text "801×9=7209" "126×2=252"
Replace-Text "435×2=870" "749×8=5992"
Replace-Text "461×2=922" "686×8=5488"
Replace-Text "993×3=2979" "763×4=3052"
Replace-Text "850×2=1700" "704×8=5632"

Replace-Text "895×4=3580" "240×5=1200"
Replace-Text "758×9=6822" "978×5=4890"
Replace-Text "305×4=1220" "136×3=408"
Replace-Text "851×2=1702" "356×3=1068"
Replace-Text "167×3=501" "874×3=2622"

Replace-Text "355×8=2840" "972×2=1944"
Replace-Text "371×5=1855" "376×9=3384"
Replace-Text "963×7=6741" "646×4=2584"
Replace-Text "134×8=1072" "173×9=1557"
Replace-Text "244×3=732" "954×8=7632"

Replace-Text "879×2=1758" "357×9=3213"
Replace-Text "410×9=3690" "998×6=5988"
Replace-Text "677×4=2708" "606×6=3636"
Replace-Text "378×6=2268" "696×6=4176"
Replace-Text "830×4=3320" "115×9=1035"
